$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were deleted from the source data (RM 232, SC 92)
$ws.Rows("28").Delete()
$ws.Rows("26").Delete()

# Apply corrected / re-imputed values for the remaining cells
$ws.Range("C2").Value = 14.9
$ws.Range("F2").Value = 18.03
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("E4").Value = -6.4
$ws.Range("E5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("C12").Value = 12.5
$ws.Range("F13").Value = 17.1
$ws.Range("C14").ClearContents()
$ws.Range("F19").ClearContents()
$ws.Range("C20").Value = 12.5
$ws.Range("C21").Value = 12.7
$ws.Range("C22").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Range("E23").Value = -7
$ws.Range("F25").Value = 16.6
$ws.Range("E27").ClearContents()
$ws.Range("F28").Value = 17.44
$ws.Range("E29").Value = -6.8
$ws.Range("B30").Value = -19.7
$ws.Range("C31").Value = 15.3
$ws.Range("F31").ClearContents()
$ws.Range("B32").ClearContents()
$ws.Range("F32").Value = 17.39
$ws.Range("C33").Value = 10.4
